$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")
$ws.Range("D15").Value = "Hibor"
